$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.934.33"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.546.09"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.49"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.35"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.19%  "
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.79"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "2.934.87"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "2.576.81"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("E16").Value = "  +7.23%  "
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "42.950.82"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.20"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.52%  "
$ws.Range("D20").Value = "0.0₃0992"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.57"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.80"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.07"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.76"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.33"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +10.06%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.74"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.21"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.03"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.85"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.37%  "
$ws.Range("E37").Value = "  -4.57%  "
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.120"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.17"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.03%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.46"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.91"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.08"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -7.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0306"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("D45").Value = "2.072.08"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.43"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.04"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").Value = "2.790.76"
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.59"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.87%  "
